# SSU last update, and inicial version of testing
#
# Applies:
#  1. Text edits / paragraph merges & deletions in the body of the
#     document (done first, while paragraph indices are still stable).
#  2. Appends a new row to the revision-history table (done last,
#     since touching the table perturbs paragraph indexing).

$d = $word.ActiveDocument

# --- Delete paragraphs "4a) ..." and "4b) ..." entirely ----------------
$p126 = $d.Paragraphs.Item(126)
$p127 = $d.Paragraphs.Item(127)
$d.Range($p126.Range.Start, $p127.Range.End).Delete()

# --- "Klikom na padajući meni ..." -> new sentence about team info ----
$d.Content.Find.Execute(
    "Klikom na padajući meni klijent dobija opciju da izabere broj kola za koji želi da vidi rezultate.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Klijent može direktno pristupiti informacijama o timu koji se nalazi u tabeli odigranih utakmica klikom na ime tog tima.",
    2) | Out-Null

# --- "... rezultati kola koje je trenutno u toku ... nijedno nije u toku." -> "... rezultati svih odigranih utakmica." ---
$d.Content.Find.Execute(
    "rezultati kola koje je trenutno u toku ili kola koje je poslednje završeno u slučaju da nijedno nije u toku.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "rezultati svih odigranih utakmica.",
    2) | Out-Null

# --- Merge the "putem sajta." paragraph with the "Odabirom kola ..." --
# --- paragraph (deleted) and the trailing tab-only paragraph ----------
$p116 = $d.Paragraphs.Item(116)
$p117 = $d.Paragraphs.Item(117)
$d.Range($p116.Range.End - 1, $p117.Range.End).Delete()
$p116b = $d.Paragraphs.Item(116)
$d.Range($p116b.Range.End - 1, $p116b.Range.End).Delete()

# --- "... putem sajta." -> "... na ovoj stranici." ---------------------
$d.Content.Find.Execute(
    "putem sajta.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "na ovoj stranici.",
    2) | Out-Null

# --- Append a new row to the revision-history table --------------------
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "6.6.2020."
$newRow.Cells.Item(2).Range.Text = "1.2"
$newRow.Cells.Item(3).Range.Text = "Ispravke u skladu sa implementacijom"
$newRow.Cells.Item(4).Range.Text = "Nikola Barjaktarević"

Write-Output "done"
